$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark wherever it currently lives.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Fix the typo: "note with an locker" -> "note with a locker"
$d.Content.Find.Execute("note with an locker", $false, $false, $false, $false,
                         $false, $true, 1, $false, "note with a locker", 2)

# 3) Re-add the "_GoBack" bookmark right after "note with a" (i.e. right
#    before " locker"), matching where Word leaves it after the last edit.
$target = $d.Content
$target.Find.Execute("note with a", $false, $false, $false, $false,
                      $false, $true, 1, $false, "", 0)
$bookmarkRange = $d.Range($target.End, $target.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
